# Applies the 'scheduled runner' numeric refresh to the per-job Leve profit
# tables (currentAveragePrice / LevePrice / LeveProfit columns H:N) across all
# eight crafting-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 855.4  # was 710.1852
$ws.Range("I28").Value = 571.95  # was 426.7
$ws.Range("J28").Value = 1989.2  # was 1520.1428
$ws.Range("K28").Value = 571.95  # was 426.7
$ws.Range("L28").Value = 1989.2  # was 1520.1428
$ws.Range("M28").Value = -86.95000000000005  # was 58.30000000000001
$ws.Range("N28").Value = -2959.2  # was -2490.1428
$ws.Range("H31").Value = 3975  # was 1234
$ws.Range("I31").Value = 3300  # was 851
$ws.Range("J31").Value = 6000  # was 2000
$ws.Range("K31").Value = 9900  # was 2553
$ws.Range("L31").Value = 18000  # was 6000
$ws.Range("M31").Value = -9670  # was -2323
$ws.Range("N31").Value = -18460  # was -6460
$ws.Range("H107").Value = 6250629  # was 940.0909
$ws.Range("I107").Value = 10000360  # was 1255.3636
$ws.Range("J107").Value = 1076  # was 624.8182
$ws.Range("K107").Value = 10000360  # was 1255.3636
$ws.Range("L107").Value = 1076  # was 624.8182
$ws.Range("M107").Value = -9998440  # was 664.6364000000001
$ws.Range("N107").Value = -4916  # was -4464.8182
$ws.Range("H111").Value = 4372.2144  # was 4828.364
$ws.Range("I111").Value = 3774.875  # was 4100
$ws.Range("J111").Value = 5168.6665  # was 5702.4
$ws.Range("K111").Value = 11324.625  # was 12300
$ws.Range("L111").Value = 15505.9995  # was 17107.2
$ws.Range("M111").Value = -8257.625  # was -9233
$ws.Range("N111").Value = -21639.9995  # was -23241.2
$ws.Range("H116").Value = 6349.88  # was 5741.8
$ws.Range("I116").Value = 7196.5  # was 4180
$ws.Range("J116").Value = 5785.467  # was 6582.769
$ws.Range("K116").Value = 7196.5  # was 4180
$ws.Range("L116").Value = 5785.467  # was 6582.769
$ws.Range("M116").Value = -3754.5  # was -738
$ws.Range("N116").Value = -12669.467  # was -13466.769
$ws.Range("H127").Value = 956.5714  # was 967.94446
$ws.Range("I127").Value = 330.25  # was 317.4
$ws.Range("J127").Value = 1791.6666  # was 1781.125
$ws.Range("K127").Value = 990.75  # was 952.1999999999999
$ws.Range("L127").Value = 5374.9998  # was 5343.375
$ws.Range("M127").Value = 3969.25  # was 4007.8
$ws.Range("N127").Value = -15294.9998  # was -15263.375
$ws.Range("H137").Value = 281627.72  # was 336603.12
$ws.Range("I137").Value = 467264.25  # was 513951.34
$ws.Range("J137").Value = 3172.9092  # was 4075.25
$ws.Range("K137").Value = 1401792.75  # was 1541854.02
$ws.Range("L137").Value = 9518.7276  # was 12225.75
$ws.Range("M137").Value = -1399242.75  # was -1539304.02
$ws.Range("N137").Value = -14618.7276  # was -17325.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2133.9  # was 2171.9656
$ws.Range("I2").Value = 1335.15  # was 1323.65
$ws.Range("J2").Value = 3731.4  # was 4057.111
$ws.Range("K2").Value = 1335.15  # was 1323.65
$ws.Range("L2").Value = 3731.4  # was 4057.111
$ws.Range("M2").Value = -1222.15  # was -1210.65
$ws.Range("N2").Value = -3957.4  # was -4283.111
$ws.Range("H74").Value = 1819.381  # was 2055.0857
$ws.Range("I74").Value = 1170  # was 1293.4814
$ws.Range("J74").Value = 4200.4443  # was 4625.5
$ws.Range("K74").Value = 1170  # was 1293.4814
$ws.Range("L74").Value = 4200.4443  # was 4625.5
$ws.Range("M74").Value = -296  # was -419.4813999999999
$ws.Range("N74").Value = -5948.4443  # was -6373.5
$ws.Range("H77").Value = 1819.381  # was 2055.0857
$ws.Range("I77").Value = 1170  # was 1293.4814
$ws.Range("J77").Value = 4200.4443  # was 4625.5
$ws.Range("K77").Value = 5850  # was 6467.406999999999
$ws.Range("L77").Value = 21002.2215  # was 23127.5
$ws.Range("M77").Value = -1482  # was -2099.406999999999
$ws.Range("N77").Value = -29738.2215  # was -31863.5
$ws.Range("H80").Value = 38030  # was 38055
$ws.Range("J80").Value = 38030  # was 38055
$ws.Range("L80").Value = 38030  # was 38055
$ws.Range("N80").Value = -40026  # was -40051
$ws.Range("H83").Value = 38030  # was 38055
$ws.Range("J83").Value = 38030  # was 38055
$ws.Range("L83").Value = 114090  # was 114165
$ws.Range("N83").Value = -124074  # was -124149
$ws.Range("H116").Value = 2133.9  # was 2171.9656
$ws.Range("I116").Value = 1335.15  # was 1323.65
$ws.Range("J116").Value = 3731.4  # was 4057.111
$ws.Range("K116").Value = 1335.15  # was 1323.65
$ws.Range("L116").Value = 3731.4  # was 4057.111
$ws.Range("M116").Value = 958.8499999999999  # was 970.3499999999999
$ws.Range("N116").Value = -8319.4  # was -8645.111000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2133.9  # was 2171.9656
$ws.Range("I3").Value = 1335.15  # was 1323.65
$ws.Range("J3").Value = 3731.4  # was 4057.111
$ws.Range("K3").Value = 1335.15  # was 1323.65
$ws.Range("L3").Value = 3731.4  # was 4057.111
$ws.Range("M3").Value = -1221.15  # was -1209.65
$ws.Range("N3").Value = -3959.4  # was -4285.111
$ws.Range("H35").Value = 30749.6  # was 30949.6
$ws.Range("J35").Value = 30749.6  # was 30949.6
$ws.Range("L35").Value = 30749.6  # was 30949.6
$ws.Range("N35").Value = -31369.6  # was -31569.6
$ws.Range("H82").Value = 27886.643  # was 28029.5
$ws.Range("J82").Value = 33817.6  # was 34017.6
$ws.Range("L82").Value = 33817.6  # was 34017.6
$ws.Range("N82").Value = -34583.6  # was -34783.6
$ws.Range("H85").Value = 27886.643  # was 28029.5
$ws.Range("J85").Value = 33817.6  # was 34017.6
$ws.Range("L85").Value = 33817.6  # was 34017.6
$ws.Range("N85").Value = -36469.6  # was -36669.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19890  # was 19895
$ws.Range("J41").Value = 23853.334  # was 23860
$ws.Range("L41").Value = 23853.334  # was 23860
$ws.Range("N41").Value = -24709.334  # was -24716
$ws.Range("H51").Value = 9082.799999999999  # was 9082.200000000001
$ws.Range("J51").Value = 9082.799999999999  # was 9082.200000000001
$ws.Range("L51").Value = 9082.799999999999  # was 9082.200000000001
$ws.Range("N51").Value = -10554.8  # was -10554.2
$ws.Range("H60").Value = 26121.46  # was 37253.285
$ws.Range("J60").Value = 26121.46  # was 37253.285
$ws.Range("L60").Value = 26121.46  # was 37253.285
$ws.Range("N60").Value = -27143.46  # was -38275.285
$ws.Range("H61").Value = 9082.799999999999  # was 9082.200000000001
$ws.Range("J61").Value = 9082.799999999999  # was 9082.200000000001
$ws.Range("L61").Value = 9082.799999999999  # was 9082.200000000001
$ws.Range("N61").Value = -9778.799999999999  # was -9778.200000000001
$ws.Range("H68").Value = 17486.334  # was 17500
$ws.Range("J68").Value = 17486.334  # was 17500
$ws.Range("L68").Value = 17486.334  # was 17500
$ws.Range("N68").Value = -18984.334  # was -18998
$ws.Range("H71").Value = 17486.334  # was 17500
$ws.Range("J71").Value = 17486.334  # was 17500
$ws.Range("L71").Value = 52459.00199999999  # was 52500
$ws.Range("N71").Value = -59947.00199999999  # was -59988
$ws.Range("H107").Value = 617.8276  # was 703.43475
$ws.Range("I107").Value = 881.2143  # was 792.2308
$ws.Range("J107").Value = 372  # was 588
$ws.Range("K107").Value = 881.2143  # was 792.2308
$ws.Range("L107").Value = 372  # was 588
$ws.Range("M107").Value = 1038.7857  # was 1127.7692
$ws.Range("N107").Value = -4212  # was -4428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 0  # was 19000
$ws.Range("J106").Value = 0  # was 19000
$ws.Range("L106").Value = 0  # was 57000
$ws.Range("N106").ClearContents()  # was -58892, now blank

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19887  # was 13954.2
$ws.Range("I57").Value = 0  # was 5055
$ws.Range("K57").Value = 0  # was 5055
$ws.Range("M57").ClearContents()  # was -4235, now blank
$ws.Range("H70").Value = 27393.666  # was 27925.342
$ws.Range("I70").Value = 33914.23  # was 34794.06
$ws.Range("K70").Value = 33914.23  # was 34794.06
$ws.Range("M70").Value = -33644.23  # was -34524.06
$ws.Range("H73").Value = 27393.666  # was 27925.342
$ws.Range("I73").Value = 33914.23  # was 34794.06
$ws.Range("K73").Value = 33914.23  # was 34794.06
$ws.Range("M73").Value = -32978.23  # was -33858.06
$ws.Range("H107").Value = 633.5625  # was 598.74194
$ws.Range("I107").Value = 441.68182  # was 464.61905
$ws.Range("J107").Value = 1055.7  # was 880.4
$ws.Range("K107").Value = 441.68182  # was 464.61905
$ws.Range("L107").Value = 1055.7  # was 880.4
$ws.Range("M107").Value = 1478.31818  # was 1455.38095
$ws.Range("N107").Value = -4895.7  # was -4720.4
$ws.Range("H123").Value = 30158.363  # was 30728.6
$ws.Range("J123").Value = 30158.363  # was 30728.6
$ws.Range("L123").Value = 30158.363  # was 30728.6
$ws.Range("N123").Value = -35058.363  # was -35628.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2607.7144  # was 1988.7241
$ws.Range("I61").Value = 1550.8  # was 1785.2106
$ws.Range("J61").Value = 5250  # was 2375.4
$ws.Range("K61").Value = 1550.8  # was 1785.2106
$ws.Range("L61").Value = 5250  # was 2375.4
$ws.Range("M61").Value = -1348.8  # was -1583.2106
$ws.Range("N61").Value = -5654  # was -2779.4
$ws.Range("H113").Value = 2607.7144  # was 1988.7241
$ws.Range("I113").Value = 1550.8  # was 1785.2106
$ws.Range("J113").Value = 5250  # was 2375.4
$ws.Range("K113").Value = 1550.8  # was 1785.2106
$ws.Range("L113").Value = 5250  # was 2375.4
$ws.Range("M113").Value = 619.2  # was 384.7893999999999
$ws.Range("N113").Value = -9590  # was -6715.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 17055.666  # was 17688.5
$ws.Range("J109").Value = 17055.666  # was 17688.5
$ws.Range("L109").Value = 17055.666  # was 17688.5
$ws.Range("N109").Value = -19829.666  # was -20462.5
$ws.Range("H113").Value = 64047.875  # was 40904.88
$ws.Range("I113").Value = 100693.3  # was 66975
$ws.Range("J113").Value = 2972.1667  # was 1799.7
$ws.Range("K113").Value = 302079.9  # was 200925
$ws.Range("L113").Value = 8916.500100000001  # was 5399.1
$ws.Range("M113").Value = -299909.9  # was -198755
$ws.Range("N113").Value = -13256.5001  # was -9739.1
